# Update regression summary table (cap_gen_year14final) with recomputed
# coefficients / stats from the re-run hourly-data model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.1034266086838253

# Row 3
$ws.Range("B3").Value = 0.1006472130890178
$ws.Range("H3").Value = 0.2040738217728431

# Row 4
$ws.Range("B4").Value = 0.09232867143567874
$ws.Range("H4").Value = 0.1957552801195041

# Row 5
$ws.Range("B5").Value = 0.05631367723607251
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 0.1597402859198978

# Row 6
$ws.Range("B6").Value = 0.05061628046209641
$ws.Range("C6").Value = 0.002775124116460936
$ws.Range("D6").Value = 7.229908160878246
$ws.Range("E6").Value = 0.01722574224865763
$ws.Range("F6").Value = 0.04516223503779269
$ws.Range("G6").Value = 0.05607032588640019
$ws.Range("H6").Value = 0.1540428891459217

# Row 7
$ws.Range("B7").Value = 0.04756014946787813
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = 0.1509867581517035

# Row 8
$ws.Range("B8").Value = 0.03583754662207645
$ws.Range("C8").Value = 0.002786738962805749
$ws.Range("D8").Value = 4.940510280842267
$ws.Range("E8").Value = 0.006472974945438711
$ws.Range("F8").Value = 0.03035778965654701
$ws.Range("G8").Value = 0.0413173035876058
$ws.Range("H8").Value = 0.1392641553059018

# Row 9
$ws.Range("B9").Value = 0.03686871014905553
$ws.Range("C9").Value = 0.001845062320910928
$ws.Range("D9").Value = 5.227821554341242
$ws.Range("E9").Value = 0.006808002356149154
$ws.Range("F9").Value = 0.03324792020700708
$ws.Range("G9").Value = 0.04048950009110399
$ws.Range("H9").Value = 0.1402953188328809

# Row 10
$ws.Range("B10").Value = 0.03901001682294634
$ws.Range("C10").Value = 0.00255422690986017
$ws.Range("D10").Value = 5.773535046456106
$ws.Range("E10").Value = 0.007460108337300352
$ws.Range("F10").Value = 0.033988306969606
$ws.Range("G10").Value = 0.04403172667628798
$ws.Range("H10").Value = 0.1424366255067717

# Row 11
$ws.Range("B11").Value = 0.03088017471332317
$ws.Range("H11").Value = 0.1343067833971485

# Row 12
$ws.Range("B12").Value = 0.05147464369771785
$ws.Range("H12").Value = 0.1549012523815432

# Row 13
$ws.Range("B13").Value = 0.0677617951436638
$ws.Range("H13").Value = 0.1711884038274891

# Row 14
$ws.Range("B14").Value = 0.07175102020736146
$ws.Range("H14").Value = 0.1751776288911868

# Row 15
$ws.Range("B15").Value = 0.07907745582211108
$ws.Range("H15").Value = 0.1825040645059364

# Row 16
$ws.Range("B16").Value = 0.0837079406595235
$ws.Range("H16").Value = 0.1871345493433488

# Row 17
$ws.Range("B17").Value = 0.08589609529377445
$ws.Range("H17").Value = 0.1893227039775998

# Row 18
$ws.Range("B18").Value = -0.1034266086838253

# Row 19
$ws.Range("B19").Value = 0.08815766389767526
$ws.Range("H19").Value = 0.1915842725815006

# Row 20
$ws.Range("B20").Value = 0.09190117280910687
$ws.Range("H20").Value = 0.1953277814929322

# Row 21
$ws.Range("B21").Value = 0.09583652875897439
$ws.Range("H21").Value = 0.1992631374427997

# Row 22
$ws.Range("B22").Value = 0.1000847894830573
$ws.Range("C22").Value = 0.007506695754597443
$ws.Range("D22").Value = 24.23474978223476
$ws.Range("E22").Value = 0.04219457408128865
$ws.Range("F22").Value = 0.08530715130748218
$ws.Range("G22").Value = 0.1148624276586325
$ws.Range("H22").Value = 0.2035113981668826

# Row 23
$ws.Range("B23").Value = 0.1036830179220661
$ws.Range("C23").Value = 0.007398999447646707
$ws.Range("D23").Value = 25.16679227113135
$ws.Range("E23").Value = 0.03962144289720126
$ws.Range("F23").Value = 0.08913937494315255
$ws.Range("G23").Value = 0.1182266609009801
$ws.Range("H23").Value = 0.2071096266058914

# Row 24
$ws.Range("B24").Value = 0.1091933794225123
$ws.Range("C24").Value = 0.007456009891838536
$ws.Range("D24").Value = 26.78362469229936
$ws.Range("E24").Value = 0.04507251405913407
$ws.Range("F24").Value = 0.09452914025756412
$ws.Range("G24").Value = 0.1238576185874605
$ws.Range("H24").Value = 0.2126199881063376

# Row 25
$ws.Range("B25").Value = 0.110467683033059
$ws.Range("C25").Value = 0.008766470598799106
$ws.Range("D25").Value = 25.59246074959471
$ws.Range("E25").Value = 0.04161387311213175
$ws.Range("F25").Value = 0.09320792228469538
$ws.Range("G25").Value = 0.1277274437814231
$ws.Range("H25").Value = 0.2138942917168843

# Row 26
$ws.Range("B26").Value = 0.111458949772634
$ws.Range("C26").Value = 0.008365031793420543
$ws.Range("D26").Value = 25.75841447348996
$ws.Range("E26").Value = 0.03879307820635577
$ws.Range("F26").Value = 0.09501340980857417
$ws.Range("G26").Value = 0.1279044897366938
$ws.Range("H26").Value = 0.2148855584564593

# Row 27
$ws.Range("B27").Value = 0.1118534318178833
$ws.Range("C27").Value = 0.007433828465952445
$ws.Range("D27").Value = 25.10410679833312
$ws.Range("E27").Value = 0.05196740885247918
$ws.Range("F27").Value = 0.09724637480122314
$ws.Range("G27").Value = 0.1264604888345434
$ws.Range("H27").Value = 0.2152800405017086

# Row 28
$ws.Range("B28").Value = 0.108244377125806
$ws.Range("C28").Value = 0.007414605910322103
$ws.Range("D28").Value = 24.00407113535745
$ws.Range("E28").Value = 0.08279491930408636
$ws.Range("F28").Value = 0.09367797055066414
$ws.Range("G28").Value = 0.1228107837009481
$ws.Range("H28").Value = 0.2116709858096313

# Row 29
$ws.Range("B29").Value = 0.04238876528774297
$ws.Range("C29").Value = 0.002351870200783593
$ws.Range("D29").Value = 6.994076770040634
$ws.Range("E29").Value = 0.006103438735143013
$ws.Range("F29").Value = 0.03776587826790539
$ws.Range("G29").Value = 0.04701165230758043
$ws.Range("H29").Value = 0.1458153739715683

